$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.086.18'
$ws.Range("D3").Value = '2.961.56'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.27'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.546'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.85%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.67'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0853'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.41'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.57%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.419.00'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  +5.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.14'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +70.97%  '
$ws.Range("D17").Value = '2.956.42'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("E18").Value = '  +2.67%  '
$ws.Range("D19").Value = '51.164.88'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.10'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.43'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").Value = '0.0₃0963'
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("E23").Value = '  +16.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.12'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.79'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.91'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.07'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -10.44%  '
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.49'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.24%  '
$ws.Range("E33").Value = '  +7.36%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.28'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '34.39'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0437'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.09%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.30'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +10.39%  '
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("E40").Value = '  +1.64%  '
$ws.Range("E41").Value = '  +3.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.51'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '124.75'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.20%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.56'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +10.38%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.59'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.273'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.38'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.00%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.02'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.063.13'
$ws.Range("E49").Value = '  +3.17%  '
$ws.Range("E50").Value = '  -8.49%  '
$ws.Range("E51").Value = '  +6.78%  '
